$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the date number format to H1 first, then propagate the same style
# (as opposed to a fresh xf per cell) to H2:H3 via copy/paste-format.
$ws.Range("H1").NumberFormat = "mm-dd-yy"
$ws.Range("H1").Copy()
$ws.Range("H2:H3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new date values (Excel serial date numbers).
$ws.Range("H1").Value = 41733
$ws.Range("H2").Value = 42129

# Move / update the active selection like the author did.
$ws.Range("H10").Select()
